# Update the "p" (significance) column values in the Mantel correlogram
# table to reflect the new x-axis distance-class label.
$d = $word.ActiveDocument

$replacements = @(
    @{ Old = "0.247"; New = "0.255" },
    @{ Old = "0.058";  New = "0.054" },
    @{ Old = "0.494"; New = "0.509" },
    @{ Old = "0.237"; New = "0.213" },
    @{ Old = "0.2";   New = "0.192" },
    @{ Old = "0.773"; New = "0.777" }
)

foreach ($r in $replacements) {
    $find = $d.Content.Find
    $find.ClearFormatting()
    $find.Execute($r.Old, $true, $true, $false, $false, $false, $true, 1, $false, $r.New, 2)
}
